$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Paragraph 7 ("Threw in a nifty background ...") -> split into several runs
# describing the new work (backgrounds per-section, text color adjustments).
# ---------------------------------------------------------------------------
$p7 = $d.Paragraphs(7)
$rng7 = $d.Range($p7.Range.Start, $p7.Range.End)
$rng7.Text = ""
$p7 = $d.Paragraphs(7)
$xmlFrag7 = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:tab/><w:t>Threw in a</w:t></w:r><w:r><w:t xml:space="preserve"> few</w:t></w:r><w:r><w:t xml:space="preserve"> nifty background</w:t></w:r><w:r><w:t>s</w:t></w:r><w:r><w:t xml:space="preserve"> I found online</w:t></w:r><w:r><w:t xml:space="preserve"> onto each of the sections</w:t></w:r><w:r><w:t xml:space="preserve">. </w:t></w:r><w:r><w:t>Adjusted text color and a few other minor things to make it work out.</w:t></w:r></w:p>
"@
$p7.Range.InsertXML($xmlFrag7)

# ---------------------------------------------------------------------------
# Paragraph 9 ("Gave my titles ...") -> fix "widthto" into "width" + " " +
# "to" as three runs, and drop the _GoBack bookmark (it will be re-added on
# the new "Put all images..." paragraph further down).
# ---------------------------------------------------------------------------
$p9 = $d.Paragraphs(9)
$rng9 = $d.Range($p9.Range.Start, $p9.Range.End)
$rng9.Text = ""
$p9 = $d.Paragraphs(9)
$xmlFrag9 = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:tab/><w:t>Gave my titles an image for background, reduced their width</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">to the width of the image and used “margin: 0 auto” to center them. </w:t></w:r></w:p>
"@
$p9.Range.InsertXML($xmlFrag9)

# ---------------------------------------------------------------------------
# New trailing paragraphs describing the remaining work:
#   (blank)
#   Put some padding on the body ...
#   (blank)
#   Put all images in an images file.   <- keeps the _GoBack bookmark
#   (blank)
#   (tab only, trailing empty paragraph)
# ---------------------------------------------------------------------------
$p9 = $d.Paragraphs(9)
$p9.Range.InsertParagraphAfter()

$p10 = $d.Paragraphs(10)
$p10.Range.InsertParagraphAfter()

$p11 = $d.Paragraphs(11)
$xmlFrag11 = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:tab/><w:t>Put some padding on the body to space things out a bit more, reduced the width to 80% for more of a margin and centered using auto.</w:t></w:r></w:p>
"@
$p11.Range.InsertXML($xmlFrag11)

$p11 = $d.Paragraphs(11)
$p11.Range.InsertParagraphAfter()

$p12 = $d.Paragraphs(12)
$p12.Range.InsertParagraphAfter()

$p13 = $d.Paragraphs(13)
$xmlFrag13 = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:tab/><w:t>Put all images in an images file.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
"@
$p13.Range.InsertXML($xmlFrag13)

$p13 = $d.Paragraphs(13)
$p13.Range.InsertParagraphAfter()

$p14 = $d.Paragraphs(14)
$p14.Range.InsertParagraphAfter()

$p15 = $d.Paragraphs(15)
$xmlFrag15 = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:tab/></w:r></w:p>
"@
$p15.Range.InsertXML($xmlFrag15)

Write-Output "Final paragraph count: $($d.Paragraphs.Count)"
foreach ($p in $d.Paragraphs) {
    Write-Output ("PARA: [" + $p.Range.Text + "]")
}
